$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '37.916.17'
$ws.Range('E2').Value = '  -0.74%  '
Set-TextValue 'D3' '2.039.00'
$ws.Range('E3').Value = '  -1.08%  '
Set-TextValue 'D5' '227.45'
$ws.Range('E5').Value = '  -1.26%  '
Set-TextValue 'D6' '0.614'
$ws.Range('E6').Value = '  -0.22%  '
Set-TextValue 'D7' '60.29'
$ws.Range('E7').Value = '  +3.63%  '
$ws.Range('E8').Value = '  -0.02%  '
Set-TextValue 'D9' '0.387'
$ws.Range('E9').Value = '  -0.35%  '
Set-TextValue 'D10' '0.0817'
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('E11').Value = '  -0.10%  '
Set-TextValue 'D12' '14.69'
$ws.Range('E12').Value = '  +0.35%  '
Set-TextValue 'D13' '2.340.92'
$ws.Range('E13').Value = '  -1.09%  '
Set-TextValue 'D14' '21.11'
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('E15').Value = '  +0.65%  '
Set-TextValue 'D16' '5.21'
$ws.Range('E16').Value = '  -1.63%  '
Set-TextValue 'D17' '2.041.30'
$ws.Range('E17').Value = '  -0.79%  '
Set-TextValue 'D18' '37.810.31'
$ws.Range('E18').Value = '  -0.68%  '
Set-TextValue 'D19' '6.09'
$ws.Range('E19').Value = '  -1.38%  '
Set-TextValue 'D20' '69.82'
$ws.Range('E20').Value = '  -0.17%  '
Set-TextValue 'D21' '0.0₃0825'
$ws.Range('E21').Value = '  -0.87%  '
Set-TextValue 'D22' '225.16'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -2.22%  '
Set-TextValue 'D25' '2.21'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('E26').Value = '  -0.50%  '
Set-TextValue 'D27' '165.29'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('E28').Value = '  -3.62%  '
Set-TextValue 'D29' '18.96'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('E30').Value = '  -6.16%  '
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('E33').Value = '  +4.04%  '
Set-TextValue 'D34' '4.50'
$ws.Range('E34').Value = '  -2.48%  '
Set-TextValue 'D35' '0.0602'
$ws.Range('E35').Value = '  -2.40%  '
Set-TextValue 'D36' '6.40'
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('E37').Value = '  -5.29%  '
Set-TextValue 'D38' '3.25'
$ws.Range('E38').Value = '  -2.75%  '
$ws.Range('E39').Value = '  -0.17%  '
Set-TextValue 'D40' '1.544.95'
$ws.Range('E40').Value = '  +4.02%  '
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '16.99'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D43' '97.22'
$ws.Range('E43').Value = '  -1.37%  '
$ws.Range('E44').Value = '  -0.93%  '
Set-TextValue 'D45' '0.0924'
$ws.Range('E45').Value = '  -2.17%  '
$ws.Range('E46').Value = '  -1.40%  '
Set-TextValue 'D47' '3.91'
$ws.Range('E47').Value = '  -5.03%  '
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('E49').Value = '  +0.06%  '
Set-TextValue 'D50' '7.12'
$ws.Range('E50').Value = '  +0.10%  '
Set-TextValue 'D51' '2.228.35'
$ws.Range('E51').Value = '  -1.12%  '
